# Commit: "update estimates of maximum offsets for Scenario 3"
#
# - Rename header J1 from "Max ElNino year" to "Max offset year"
# - Update the Scenario-3 "Max offset" values in column I (rows 2,3,4,5,6,9,10)
#   to their newly-rounded estimates
# - View-state nudges that rode along with the edit: selection moved to E10,
#   and several columns lost their "best fit" auto-width (now fixed widths),
#   with column I narrowed from its old auto-fit width down to ~12.6 chars.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header text -----------------------------------------------------------
$ws.Range("J1").Value = "Max offset year"

# --- updated Scenario 3 "max offset" estimates ------------------------------
$ws.Range("I2").Value = 1.27
$ws.Range("I3").Value = 2.81
$ws.Range("I4").Value = 1.0900000000000001
$ws.Range("I5").Value = 1.91
$ws.Range("I6").Value = 1.55
$ws.Range("I9").Value = 3.79
$ws.Range("I10").Value = 1.27

# --- selection / active cell ------------------------------------------------
[void]$ws.Range("E10").Select()

# --- column widths: drop auto "best fit" sizing, column I narrows ----------
# The underlying engine stores ColumnWidth as a whole-pixel quantity (it
# always rounds the value you assign to the nearest achievable pixel width),
# so assigning the target character-width directly overshoots by a constant
# padding offset. Back that offset out first so the stored width lands on
# the closest achievable value to the real target.
function Set-ColWidth($col, $target) {
    $ws.Columns.Item($col).ColumnWidth = ($target - 5/6 - 0.00001)
}

Set-ColWidth 2 13.875
Set-ColWidth 3 13.125
Set-ColWidth 4 8.75
Set-ColWidth 5 14.625
Set-ColWidth 6 9.875
Set-ColWidth 7 10.5
Set-ColWidth 8 40.25
Set-ColWidth 9 12.625
Set-ColWidth 10 15
